# New version of data feeds: update the RSE sheet's industry-balance
# factors and make RSE the active sheet/selection (matches the author's
# commit touching "RSE" values and workbook view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSE")

# Update values in columns B and C for the "data feeds" rows.
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 5
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 5

# Make RSE the active sheet and select B9 (matches the saved view state).
$ws.Activate()
$ws.Range("B9").Select()
